$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = @(0.003994804209775715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1, 1.619054046815646)
    3 = @(1.459612070389937, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0, 2.156321060280252)
    4 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    5 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    6 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 8.660232485948974, 1, 15.68806981981553)
    7 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 1, 29.84159230404497)
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colIndex = 2 + $i  # B=2, C=3, D=4, E=5, F=6, G=7
        $ws.Cells.Item($row, $colIndex).Value = $cols[$i]
    }
}
